# -----------------------------------------------------------------------
# "add today's qubit results" -- append the 2019-11-01 Qubit run
# (24 samples, run_ID "2019-11-01_172026") to the bottom of the
# qubit-iso data table (previously ending at row 494).
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 495
$lastRow  = 518

$rows = @(
    [pscustomobject]@{ Row=495; TestName='Sample_#191101-172409'; TestDate=43770.725104166668; Conc=405; OrigConc=40.5; UlUsed=264 },
    [pscustomobject]@{ Row=496; TestName='Sample_#191101-172359'; TestDate=43770.724988425929; Conc=166; OrigConc=16.600000000000001; UlUsed=368 },
    [pscustomobject]@{ Row=497; TestName='Sample_#191101-172350'; TestDate=43770.72488425926; Conc=271; OrigConc=27.1; UlUsed=289 },
    [pscustomobject]@{ Row=498; TestName='Sample_#191101-172340'; TestDate=43770.724768518521; Conc=361; OrigConc=36.1; UlUsed=273 },
    [pscustomobject]@{ Row=499; TestName='Sample_#191101-172331'; TestDate=43770.724664351852; Conc=207; OrigConc=20.7; UlUsed=363 },
    [pscustomobject]@{ Row=500; TestName='Sample_#191101-172323'; TestDate=43770.72457175926; Conc=190; OrigConc=19; UlUsed=283 },
    [pscustomobject]@{ Row=501; TestName='Sample_#191101-172314'; TestDate=43770.72446759259; Conc=370; OrigConc=37; UlUsed=365 },
    [pscustomobject]@{ Row=502; TestName='Sample_#191101-172306'; TestDate=43770.724374999998; Conc=232; OrigConc=23.2; UlUsed=375 },
    [pscustomobject]@{ Row=503; TestName='Sample_#191101-172258'; TestDate=43770.724282407406; Conc=409; OrigConc=40.9; UlUsed=297 },
    [pscustomobject]@{ Row=504; TestName='Sample_#191101-172250'; TestDate=43770.724189814813; Conc=266; OrigConc=26.6; UlUsed=263 },
    [pscustomobject]@{ Row=505; TestName='Sample_#191101-172242'; TestDate=43770.724097222221; Conc=285; OrigConc=28.5; UlUsed=371 },
    [pscustomobject]@{ Row=506; TestName='Sample_#191101-172234'; TestDate=43770.724004629628; Conc=510; OrigConc=51; UlUsed=286 },
    [pscustomobject]@{ Row=507; TestName='Sample_#191101-172225'; TestDate=43770.723900462966; Conc=177; OrigConc=17.7; UlUsed=278 },
    [pscustomobject]@{ Row=508; TestName='Sample_#191101-172216'; TestDate=43770.723796296297; Conc=284; OrigConc=28.4; UlUsed=254 },
    [pscustomobject]@{ Row=509; TestName='Sample_#191101-172207'; TestDate=43770.723692129628; Conc=420; OrigConc=42; UlUsed=203 },
    [pscustomobject]@{ Row=510; TestName='Sample_#191101-172158'; TestDate=43770.723587962966; Conc=444; OrigConc=44.4; UlUsed=216 },
    [pscustomobject]@{ Row=511; TestName='Sample_#191101-172149'; TestDate=43770.723483796297; Conc=315; OrigConc=31.5; UlUsed=245 },
    [pscustomobject]@{ Row=512; TestName='Sample_#191101-172140'; TestDate=43770.723379629628; Conc=203; OrigConc=20.3; UlUsed=239 },
    [pscustomobject]@{ Row=513; TestName='Sample_#191101-172131'; TestDate=43770.723275462966; Conc=125; OrigConc=12.5; UlUsed=218 },
    [pscustomobject]@{ Row=514; TestName='Sample_#191101-172123'; TestDate=43770.723182870373; Conc=316; OrigConc=31.6; UlUsed=250 },
    [pscustomobject]@{ Row=515; TestName='Sample_#191101-172113'; TestDate=43770.723067129627; Conc=287; OrigConc=28.7; UlUsed=316 },
    [pscustomobject]@{ Row=516; TestName='Sample_#191101-172105'; TestDate=43770.722974537035; Conc=272; OrigConc=27.2; UlUsed=257 },
    [pscustomobject]@{ Row=517; TestName='Sample_#191101-172057'; TestDate=43770.722881944443; Conc=280; OrigConc=28; UlUsed=228 },
    [pscustomobject]@{ Row=518; TestName='Sample_#191101-172048'; TestDate=43770.722777777781; Conc=129; OrigConc=12.9; UlUsed=220 }
)


# --- 1. clone formatting from the last existing row (494) ----------------
# so the new rows land on the same cellXfs entries already used by the
# sheet (date format on column D, the Arial-10 style on I/K/L/M) instead
# of Excel fabricating brand-new styles.
$ws.Range("A494:N494").Copy() | Out-Null
$ws.Range("A495:N518").PasteSpecial(-4122) | Out-Null
$ws.Range("R494").Copy() | Out-Null
$ws.Range("R495:R518").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- 2. column A (run_ID) first, so the new shared-string entry for the
# run timestamp is interned before the per-sample test names ------------
foreach ($row in $rows) {
    $ws.Cells.Item($row.Row, 1).Value = "2019-11-01_172026"
}

# --- 3. column C (test_name / sample id) in row order -------------------
foreach ($row in $rows) {
    $ws.Cells.Item($row.Row, 3).Value = $row.TestName
}

# --- 4. remaining columns -------------------------------------------------
foreach ($row in $rows) {
    $r = $row.Row

    $ws.Cells.Item($r, 2).Value  = "RNA High sensitivity"        # assay_name
    $ws.Cells.Item($r, 4).Value  = $row.TestDate                 # test_date
    $ws.Cells.Item($r, 4).NumberFormat = "m/d/yy h:mm"
    $ws.Cells.Item($r, 5).Value  = $row.Conc                     # qubit_tube_conc_ng.ml
    $ws.Cells.Item($r, 6).Value  = $row.OrigConc                 # original_sample_conc_ng.ul
    $ws.Cells.Item($r, 7).Value  = 2                              # sample_vol_ul
    $ws.Cells.Item($r, 8).Value  = 100                            # dilution_factor
    $ws.Cells.Item($r, 9).Value  = $row.UlUsed                   # ul_sample-used
    $ws.Cells.Item($r, 10).Value = "Zymo_microprep"               # extraction_method
    $ws.Cells.Item($r, 11).Value = 35                             # elution_vol_ul
    $ws.Cells.Item($r, 12).Value = 15                             # (L) incubation/other constant
    $ws.Cells.Item($r, 13).Formula = "=(F$r)*(L$r-G$r)"           # total-yield_ng
    $ws.Cells.Item($r, 14).Value = "pellet"                       # sample_type
}

# --- 5. column R (notebook post link) last, after every C is in place ---
foreach ($row in $rows) {
    $ws.Cells.Item($row.Row, 18).Value = "https://grace-ac.github.io/rna-extractions-day12-qubitresults/"
}

# --- 6. move the active selection to mirror the author's final cursor ---
$ws.Range("P511").Select()
